# Generate Report for Handback
# Updates the localization-status report: marks rows as handed back
# (in sync with en-US), refreshes the handback timestamps, clears the
# stale "handback not latest" error details, and widens the
# date/error-detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 29.1666666666667
$ws.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws.Columns.Item(16).ColumnWidth = 12.8333333333333
$ws.Range("K2").Value = "2016-08-19 06:47:49"
$ws.Range("P2").Value = ""

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws.Columns.Item(16).ColumnWidth = 12.8333333333333
$ws.Range("K2").Value = "2016-08-19 06:47:56"
$ws.Range("P2").Value = ""

# --- Status text, shared across Overview/zh-cn/de-de sheets -----------
$wb.Worksheets.Item("Overview").Range("E2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("Overview").Range("F2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("zh-cn").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C2").Value = "Handed back: in sync with en-US"
